$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "Palavras"
$ws.Range("A2").Value = "Algo"
$ws.Range("A3").Value = "Maio"
